$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in the new "request" row (row 5) and other newly-populated cells ---
# The order of these assignments controls the order new shared strings are
# appended in sharedStrings.xml, so it mirrors the target workbook's layout.

# C4: status of the first record becomes "בבדיקה" (under review)
$ws.Range("C4").Value = "בבדיקה"

# C5: status of the new record
$ws.Range("C5").Value = "הוגשו חוו""ד"

# J5: approved yes/no column
$ws.Range("J5").Value = "כן"

# M5: short project description for the new record
$ws.Range("M5").Value = "אפליקציה לניטור פסיבי והתרעה מוקדמת על הופעת סיכון לשבץ מוחי"

# M6: short project description for an existing record
$ws.Range("M6").Value = "ניטור חדש של האינטרנט"

# J8: approved yes/no column for an existing record
$ws.Range("J8").Value = "לא"

# K3: header text, trailing space removed ("גובה מענק שאושר ") -> ("גובה מענק שאושר")
$ws.Range("K3").Value = "גובה מענק שאושר"

# I5: date the office reviewer's opinion was submitted
$ws.Range("I5").Value = "22/3/15"

# K5: approved grant amount (number)
$ws.Range("K5").Value = 10000

# L5: grant percentage, formatted as a percentage
$ws.Range("L5").Value = 0.5
$ws.Range("L5").NumberFormat = "0%"

# --- Column / row sizing & selection tweaks ---
$ws.Columns.Item(11).ColumnWidth = 13.14
$ws.Rows.Item(5).RowHeight = 30.75

[void]$ws.Range("I13").Select()
